# DEV-2004 Simplify Import Tasks - Remove Logstash
# Update XLSX sample: renumber the sample SKU rows from 1-4 to 6-9
# (sku/name/barcode for rows 2-5) and move the saved cursor/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: test-sku-1 / Test product 1 / BARCODE-0001 -> ...-6
$ws.Range("A2").Value = "test-sku-6"
$ws.Range("B2").Value = "Test product 6"
$ws.Range("C2").Value = "BARCODE-0006"

# Row 3: test-sku-2 / Test product 2 / BARCODE-0002 -> ...-7
$ws.Range("A3").Value = "test-sku-7"
$ws.Range("B3").Value = "Test product 7"
$ws.Range("C3").Value = "BARCODE-0007"

# Row 4: test-sku-3 / Test product 3 / BARCODE-0003 -> ...-8
$ws.Range("A4").Value = "test-sku-8"
$ws.Range("B4").Value = "Test product 8"
$ws.Range("C4").Value = "BARCODE-0008"

# Row 5: test-sku-4 / Test product 4 / BARCODE-0004 -> ...-9
$ws.Range("A5").Value = "test-sku-9"
$ws.Range("B5").Value = "Test product 9"
$ws.Range("C5").Value = "BARCODE-0009"

# Restore the saved cursor position / selection (was R6, now K27).
[void]$ws.Range("K27").Select()
